$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "http://www.avclub.com/1798447330"
$ws.Range("B3").Value = 364
$ws.Range("C3").Value = 283
$ws.Range("D3").Value = 66
$ws.Range("E3").Value = 217
$ws.Range("F3").Value = 3.439393939393939
$ws.Range("G3").Value = 54.43939393939394
$ws.Range("H3").Value = 303.0454545454546
$ws.Range("I3").Value = 2.889400921658986
$ws.Range("J3").Value = 42.66820276497696
$ws.Range("K3").Value = 235.9723502304147
$ws.Range("L3").Value = 2

# Row 4
$ws.Range("A4").Value = "http://www.avclub.com/1798505721"
$ws.Range("B4").Value = 289
$ws.Range("C4").Value = 229
$ws.Range("D4").Value = 63
$ws.Range("E4").Value = 166
$ws.Range("F4").Value = 5.714285714285714
$ws.Range("G4").Value = 82.73015873015873
$ws.Range("H4").Value = 467.5873015873016
$ws.Range("I4").Value = 2.801204819277109
$ws.Range("J4").Value = 62.47590361445783
$ws.Range("K4").Value = 352.0602409638554
$ws.Range("L4").Value = 0

# Row 5
$ws.Range("A5").Value = "http://www.avclub.com/1798517837"
$ws.Range("B5").Value = 88
$ws.Range("C5").Value = 84
$ws.Range("D5").Value = 17
$ws.Range("E5").Value = 67
$ws.Range("F5").Value = 2.294117647058823
$ws.Range("G5").Value = 23
$ws.Range("H5").Value = 130.5294117647059
$ws.Range("I5").Value = 2.238805970149254
$ws.Range("J5").Value = 25.71641791044776
$ws.Range("K5").Value = 142.0298507462687
$ws.Range("L5").Value = 3

# Hyperlinks
$ws.Hyperlinks.Add($ws.Range("A3"), "http://www.avclub.com/1798447330")
$ws.Hyperlinks.Add($ws.Range("A4"), "http://www.avclub.com/1798505721")
$ws.Hyperlinks.Add($ws.Range("A5"), "http://www.avclub.com/1798517837")

# Hyperlinks.Add applies the built-in "Hyperlink" cell style (blue/underline);
# the source data keeps the default styling, so strip it back off.
$ws.Range("A3").ClearFormats()
$ws.Range("A4").ClearFormats()
$ws.Range("A5").ClearFormats()
